$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'245.94"
$ws.Range("D3").Value = "'24.19"
$ws.Range("D4").Value = "'5.357"
$ws.Range("D5").Value = "'0.05744"
$ws.Range("D6").Value = "'6.473"
$ws.Range("D7").Value = "'3.124"
$ws.Range("D8").Value = "'0.8177"
$ws.Range("D9").Value = "'0.8706"
$ws.Range("D11").Value = "'0.06978"
$ws.Range("D13").Value = "'0.02915"
$ws.Range("D14").Value = "'0.09408"
$ws.Range("D15").Value = "'3.730"
$ws.Range("D16").Value = "'0.001558"
$ws.Range("D17").Value = "'0.04696"
$ws.Range("D18").Value = "'0.0006004"
$ws.Range("E18").Value = "17OneONEWorstin24h"
$ws.Range("D19").Value = "'0.006213"
$ws.Range("D21").Value = "'0.003911"
$ws.Range("D22").Value = "'0.00008790"
$ws.Range("D25").Value = "'0.3164"
$ws.Range("D26").Value = "'0.1312"
$ws.Range("D28").Value = "'0.0003011"
$ws.Range("D40").Value = "'0.03717"
$ws.Range("D41").Value = "'0.006403"
$ws.Range("D42").Value = "'0.1058"
$ws.Range("D43").Value = "'0.002997"
$ws.Range("D44").Value = "'0.008697"
$ws.Range("D45").Value = "'0.00005263"
$ws.Range("D47").Value = "'0.3896"
$ws.Range("D48").Value = "'0.003022"
$ws.Range("E48").Value = "47BOLOBOLO"
$ws.Range("D49").Value = "'0.00002098"
$ws.Range("D50").Value = "'0.0001998"
